$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 292730
$ws.Range("AA2").Value = "'2007-11-09"
$ws.Range("AA2").Style = "Normal"
$ws.Range("AH2").Value = 'Granskog'
$ws.Range("AI2").Value = 'Gransumpskog'
$ws.Range("AX2").Value = 'Hans Sundström'
$ws.Range("B2").Value = 79433
$ws.Range("E2").Value = 1049
$ws.Range("F2").Value = 'Kortskaftad ärgspik'
$ws.Range("G2").Value = 'Microcalicium ahlneri'
$ws.Range("H2").Value = 'Tibell'
$ws.Range("P2").Value = 'Täljeån, Mpd'
$ws.Range("Q2").Value = 540844.6604352774
$ws.Range("R2").Value = 6941718.921423005
$ws.Range("S2").Value = 25
$ws.Range("Y2").Value = "'2007-11-09"
$ws.Range("Y2").Style = "Normal"
$ws.Range("AC2").ClearContents()
$ws.Range("AN2").ClearContents()
$ws.Range("AO2").ClearContents()
$ws.Range("AR2").ClearContents()

# Row 4
$ws.Range("A4").Value = 1901472
$ws.Range("B4").Value = 78569
$ws.Range("D4").Value = 'NT'
$ws.Range("E4").Value = 6458
$ws.Range("F4").Value = 'Lunglav'
$ws.Range("G4").Value = 'Lobaria pulmonaria'
$ws.Range("H4").Value = '(L.) Hoffm.'
$ws.Range("Q4").Value = 540812.4636330464
$ws.Range("R4").Value = 6941674.310578943
$ws.Range("AH4").ClearContents()

# Row 5
$ws.Range("A5").Value = 93500
$ws.Range("AA5").Value = "'2009-09-02"
$ws.Range("AA5").Style = "Normal"
$ws.Range("AC5").Value = 'Jonas Salmonsson'
$ws.Range("AN5").Value = 2
$ws.Range("AO5").Value = '2 substratenheter # Timmer'
$ws.Range("AR5").Value = ""
$ws.Range("AX5").Value = 'Via Andreas Karlberg'
$ws.Range("B5").Value = 94121
$ws.Range("D5").Value = 'NT'
$ws.Range("E5").Value = 53
$ws.Range("F5").Value = 'Vedtrappmossa'
$ws.Range("G5").Value = 'Crossocalyx hellerianus'
$ws.Range("H5").Value = '(Nees ex Lindenb.) Meyl.'
$ws.Range("P5").Value = 'Bäcken vid Mellantjärnsbodarna, Mpd'
$ws.Range("Q5").Value = 540775.1640602688
$ws.Range("R5").Value = 6941931.758068252
$ws.Range("S5").Value = 10
$ws.Range("Y5").Value = "'2009-09-02"
$ws.Range("Y5").Style = "Normal"
$ws.Range("AH5").ClearContents()

# Row 6
$ws.Range("A6").Value = 1866241
$ws.Range("B6").Value = 73678
$ws.Range("D6").Value = 'LC'
$ws.Range("E6").Value = 6439
$ws.Range("F6").Value = 'Gulnål'
$ws.Range("G6").Value = 'Chaenotheca brachypoda'
$ws.Range("H6").Value = '(Ach.) Tibell'
$ws.Range("Q6").Value = 540835.5837245358
$ws.Range("R6").Value = 6941669.529265426

# Row 7
$ws.Range("A7").Value = 1672419
$ws.Range("AH7").Value = 'Granskog'
$ws.Range("AI7").Value = 'Gransumpskog'
$ws.Range("AO7").Value = 'Granlåga'
$ws.Range("B7").Value = 89356
$ws.Range("D7").Value = 'LC'
$ws.Range("E7").Value = 5447
$ws.Range("F7").Value = 'Vedticka'
$ws.Range("G7").Value = 'Fuscoporia viticola'
$ws.Range("H7").Value = '(Schwein.) Murrill'
$ws.Range("Q7").Value = 540838.9380165208
$ws.Range("R7").Value = 6941696.743931145

# Row 8
$ws.Range("A8").Value = 1866240
$ws.Range("Q8").Value = 540827.0929000516
$ws.Range("R8").Value = 6941759.694429157
$ws.Range("AI8").ClearContents()

# Row 9
$ws.Range("A9").Value = 168997
$ws.Range("B9").Value = 73685
$ws.Range("D9").Value = 'VU'
$ws.Range("E9").Value = 492
$ws.Range("F9").Value = 'Smalskaftslav'
$ws.Range("G9").Value = 'Chaenotheca gracilenta'
$ws.Range("H9").Value = '(Ach.) J.Mattsson & Middelb.'
$ws.Range("Q9").Value = 540827.0929000516
$ws.Range("R9").Value = 6941759.694429157
$ws.Range("AI9").ClearContents()
$ws.Range("AO9").ClearContents()
